$p = $ppt.ActivePresentation
try {
  $hm = $p.HandoutMaster
  Write-Host "got handout master"
  $hm | Get-Member
} catch {
  Write-Host "ERR $_"
}
